# [MOD] Tabelle s1 e 2 - Completata Strategia 3
# Aggiungere seconda strategia per la somma
#
# Foglio1 (Tabelle_Strategia1.xlsx):
#  - C3 (tempo seriale per N=1000) is updated to reflect the new
#    serial-time measurement, now shown with 3 decimal digits.
#  - F5 is touched/formatted (underlined) as the starting point for the
#    new "Strategia 3" section being prepared.
#  - Dependent ratios (Sp = C3/C.. and Ep = Sp/P) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update the serial execution time for N = 1000 and show it with one more
# decimal of precision (0.000 instead of 0.00).
$ws.Range("C3").Value = 1.4983
$ws.Range("C3").NumberFormat = "0.000"

# Mark the beginning of the new "Strategia 3" table by underlining F5
# (left as an empty, formatted cell for now).
$ws.Range("F5").Font.Underline = 2

# Move the active selection to the new cell, matching where work stopped.
[void]$ws.Range("F5").Select()
